$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.003078177322033415
$ws.Range("C2").Value = 0.04103571897497393
$ws.Range("D2").Value = 0.1496068669990043
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 14.0575672340967
